$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the text must be
# pinned as Text first (so Excel does not silently reinterpret a
# numeric-looking price string like "232.68" as a Number).
$updates = @(
    @{ Cell = 'D2'; Value = '30.216.90'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -0.83%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.840.06'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -1.58%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '232.68'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.58%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.4668'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -3.26%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2717'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -3.06%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.06275'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -3.65%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '1.840.40'; ForceText = $false }
    @{ Cell = 'E10'; Value = '  -1.37%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.07415'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -0.48%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '16.08'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -1.15%  '; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -3.03%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '83.64'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -4.10%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.6193'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -3.63%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '30.147.60'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -0.98%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '225.53'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -2.52%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '0.000007280'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -2.86%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '12.33'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -5.30%  '; ForceText = $false }
    @{ Cell = 'B21'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Cell = 'D21'; Value = '2.080.63'; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -0.86%  '; ForceText = $false }
    @{ Cell = 'B22'; Value = 'BinanceUSD'; ForceText = $false }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText = $false }
    @{ Cell = 'D22'; Value = '1.000'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '4.891'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -5.04%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '5.851'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -4.17%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '9.178'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -1.83%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '164.43'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -3.55%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '17.73'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -3.41%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.860'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -2.43%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '0.1035'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -1.51%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  -0.69%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '4.069'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -4.78%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '3.805'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -4.59%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.04816'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -3.36%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.140'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -3.46%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.7071'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -4.77%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '2.703'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -0.28%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.01867'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -3.33%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '2.649'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.55%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  -2.70%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '1.914'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -6.58%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'PaxDollar'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; ForceText = $false }
    @{ Cell = 'D41'; Value = '1.002'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +0.53%  '; ForceText = $false }
    @{ Cell = 'B42'; Value = 'Quant'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false }
    @{ Cell = 'D42'; Value = '104.13'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -1.84%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '5.523'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -1.10%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.4011'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -4.44%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '7.022'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -2.82%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  -2.91%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '59.66'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -3.98%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '8.543'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -3.95%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '32.86'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -2.25%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.05514'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -2.37%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '1.356'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -4.78%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
